$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: Sending cluster FAPs -> Artn -> Gfra3 -> Target cluster ECs (new) ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.306751666666667
$ws.Range("H2").Value = 6.920255
$ws.Range("I2").Value = 0.8617934705859058
$ws.Range("J2").Value = 0.8617934705859057
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.024394
$ws.Range("N2").Value = 0.073182
$ws.Range("O2").Value = 0.01368854235971825
$ws.Range("P2").Value = 0.01368854235971825
$ws.Range("Q2").Value = 0.05627090015666666
$ws.Range("R2").Value = 0.50643810141
$ws.Range("S2").Value = 0.01179669642744378
$ws.Range("T2").Value = 0.01179669642744378

# --- Row 3: Sending cluster FAPs -> Artn -> Gfra3 -> Target cluster sCs ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.306751666666667
$ws.Range("H3").Value = 6.920255
$ws.Range("I3").Value = 0.8617934705859058
$ws.Range("J3").Value = 0.8617934705859057
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.757680333333333
$ws.Range("N3").Value = 5.273041
$ws.Range("O3").Value = 0.9863114576402817
$ws.Range("P3").Value = 0.9863114576402817
$ws.Range("Q3").Value = 4.054532038383889
$ws.Range("R3").Value = 36.490788345455
$ws.Range("S3").Value = 0.849996774158462
$ws.Range("T3").Value = 0.8499967741584619

# --- Row 4 (new): Sending cluster sCs -> Artn -> Gfra3 -> Target cluster ECs ---
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3699356666666667
$ws.Range("H4").Value = 1.109807
$ws.Range("I4").Value = 0.1382065294140942
$ws.Range("J4").Value = 0.1382065294140942
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.024394
$ws.Range("N4").Value = 0.073182
$ws.Range("O4").Value = 0.01368854235971825
$ws.Range("P4").Value = 0.01368854235971825
$ws.Range("Q4").Value = 0.009024210652666667
$ws.Range("R4").Value = 0.081217895874
$ws.Range("S4").Value = 0.001891845932274475
$ws.Range("T4").Value = 0.001891845932274475

# --- Row 5 (new): Sending cluster sCs -> Artn -> Gfra3 -> Target cluster sCs ---
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3699356666666667
$ws.Range("H5").Value = 1.109807
$ws.Range("I5").Value = 0.1382065294140942
$ws.Range("J5").Value = 0.1382065294140942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.757680333333333
$ws.Range("N5").Value = 5.273041
$ws.Range("O5").Value = 0.9863114576402817
$ws.Range("P5").Value = 0.9863114576402817
$ws.Range("Q5").Value = 0.6502286458985556
$ws.Range("R5").Value = 5.852057813087
$ws.Range("S5").Value = 0.1363146834818197
$ws.Range("T5").Value = 0.1363146834818197
